# Applies the "window" and "tab" sheet additions to TestData.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add sheet "window" right after the last existing sheet ("redbus")
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$windowSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$windowSheet.Name = "window"

# Header row (row 1) - yellow fill, like the header rows on the other sheets
$windowHeaders = @(
    "child window title",
    "parent window title",
    "child window txtfield",
    "parent window txtfield",
    "confirmation message",
    "url number"
)
for ($i = 0; $i -lt $windowHeaders.Length; $i++) {
    $cell = $windowSheet.Cells.Item(1, $i + 1)
    $cell.Value = $windowHeaders[$i]
    $cell.Interior.Color = 65535
}

# Row 2
$windowSheet.Cells.Item(2, 1).Value = "Basic Controls - H Y R Tutorials"
$windowSheet.Cells.Item(2, 2).Value = "Window Handles Practice - H Y R Tutorials"
$windowSheet.Cells.Item(2, 3).Value = "child txtxfield"
$windowSheet.Cells.Item(2, 4).Value = "parent txtfield"
$windowSheet.Cells.Item(2, 5).Value = "test case passed"
$windowSheet.Cells.Item(2, 6).Value = "url3"

# Row 3
$windowSheet.Cells.Item(3, 1).Value = "XPath Practice - H Y R Tutorials"
$windowSheet.Cells.Item(3, 3).Value = "child txtxfield2"

$windowSheet.Range("A1:F3").EntireColumn.AutoFit() | Out-Null
$windowSheet.Range("D1").Select() | Out-Null

# ---------------------------------------------------------------------
# Add sheet "tab" right after "window"
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$tabSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$tabSheet.Name = "tab"

# Header row (row 1) - yellow fill
$tabHeaders = @(
    "url number",
    "parent tab",
    "child tab",
    "messgae",
    "alert message",
    "parent txtfield"
)
for ($i = 0; $i -lt $tabHeaders.Length; $i++) {
    $cell = $tabSheet.Cells.Item(1, $i + 1)
    $cell.Value = $tabHeaders[$i]
    $cell.Interior.Color = 65535
}

# Row 2
$tabSheet.Cells.Item(2, 1).Value = "url3"
$tabSheet.Cells.Item(2, 2).Value = "Window Handles Practice - H Y R Tutorials"
$tabSheet.Cells.Item(2, 3).Value = "AlertsDemo - H Y R Tutorials"
$tabSheet.Cells.Item(2, 4).Value = "test case passed"
$tabSheet.Cells.Item(2, 5).Value = "I am an alert box!"
$tabSheet.Cells.Item(2, 6).Value = "This is parent tab"

# Row 3
$tabSheet.Cells.Item(3, 3).Value = "XPath Practice - H Y R Tutorials"

# Row 4
$tabSheet.Cells.Item(4, 3).Value = "Basic Controls - H Y R Tutorials"
$tabSheet.Cells.Item(4, 5).Value = "Hello!"

$tabSheet.Range("A1:F4").EntireColumn.AutoFit() | Out-Null
$tabSheet.Range("C7").Select() | Out-Null
